$d = $word.ActiveDocument

# --- Merge the "Versi" + "on" runs into a single "Version" run. ---
# A same-text assignment is a no-op in this engine, so nudge the text to
# something different first and then correct it; that forces the range's
# underlying runs to be rewritten as one.
$r = $d.Range(0, 7)
$r.Text = "Versionx"
$r = $d.Range(0, 8)
$r.Text = "Version"

# --- "Version 2." -> "Version 1." ---
# Replace just the " 2" run's text with " 1." (3 chars for a 2-char
# range). This merges the " 2" run into a new " 1." run while leaving the
# "_GoBack" bookmark between it and the trailing "." run untouched.
$r = $d.Range(7, 9)
$r.Text = " 1."

# The original trailing "." run (now duplicated after the insert above)
# is no longer needed; delete that character, leaving the bookmark
# sitting right after " 1." as in the target.
$r2 = $d.Range(10, 11)
$r2.Text = ""
